$wb = $excel.ActiveWorkbook

# --- Worksheet handles ---
$wsRanged = $wb.Worksheets.Item("Templar Ranged Weapons")
$wsModels = $wb.Worksheets.Item("Templar Models")

# --- Templar Models: add a "Test Ork" row for morale testing (added first so
#     its shared string lands before "Test Gun") ---
$wsModels.Range("A7").Value = "Test Ork"
$wsModels.Range("B7").Value = 5
$wsModels.Range("C7").Value = 3
$wsModels.Range("D7").Value = 5
$wsModels.Range("E7").Value = 4
$wsModels.Range("F7").Value = 4
$wsModels.Range("G7").Value = 1
$wsModels.Range("H7").Value = 2
$wsModels.Range("I7").Value = 3
$wsModels.Range("J7").Value = 6
$wsModels.Range("L7").Value = 12

# --- Templar Ranged Weapons: bump the Bolter's Shots from 1 to 2, and add a
#     new "Test Gun" weapon row for testing ---
$wsRanged.Range("E3").Value = 2

$wsRanged.Range("A13").Value = "Test Gun"
$wsRanged.Range("B13").Value = 12
$wsRanged.Range("C13").Value = "Assault"
$wsRanged.Range("D13").Value = 0
$wsRanged.Range("E13").Value = 2
$wsRanged.Range("F13").Value = 5
$wsRanged.Range("G13").Value = 2
$wsRanged.Range("H13").Value = 0
$wsRanged.Range("I13").Value = 1

# --- Selection / active-tab bookkeeping to match the saved UI state ---
$wsModels.Range("J7").Select()
$wsRanged.Range("E16").Select()
